$wb = $excel.ActiveWorkbook

# Turn on iterative calculation (workbook now allows circular refs to settle).
$excel.Iteration = $true
$excel.MaxChange = 0.00001

# ---------------------------------------------------------------------------
# "manufacture" sheet: update/author the notes, refresh the GREET-derived
# electricity manufacturing-energy number, and add the scratch-work formulas
# the author used to sanity check it.
# ---------------------------------------------------------------------------
$wsManufacture = $wb.Worksheets.Item("manufacture")

# Note ordering matters for shared-string allocation: the electric-tab note
# was (re)written before the petroleum-tab note.
$wsManufacture.Range("A3").Value = "line 541 from 'electric' tab, column G (CA specific value; 1,754,558)"
$wsManufacture.Range("A2").Value = "line 263 from 'petroleum' tab (CA feedstock (44,763) + CA fuel (diesel > gas) "
$wsManufacture.Range("A4").Value = "also looked at the eqn in Results tab for electricity AN14, they add fuel and feedstock together, but subtract 1,000,000 and get 1,072,402? Used that for now"

# Revised electricity manufacturing-energy value (was the old placeholder
# 1751558.4793751165 pulled from electric-by-state!A11).
$wsManufacture.Range("D9").Value = 1072402

# Scratch-work formulas kept below the table for reference.
$wsManufacture.Range("D18").Formula = "=1*3.3*0.23"
$wsManufacture.Range("D19").Formula = "=1*4.2*0.21"
$wsManufacture.Range("D20").Formula = "=1.1*1.07"

$wsManufacture.Range("D21").Select()

# ---------------------------------------------------------------------------
# Final view state when the workbook was saved: "conversion-eff" tab active
# with E7 selected (and "energy" no longer the active tab).
# ---------------------------------------------------------------------------
$wsConversion = $wb.Worksheets.Item("conversion-eff")
$wsConversion.Range("E7").Select()
